# "Generate Report for Handoff" — refresh the handoff-status report.
#
# The localization report had previously recorded a handback state
# ("Handed back: in sync with en-US") together with the timestamps of the
# last handback/handoff run. Regenerating the report for a fresh handoff
# flips the status back to "Ready for handoff" and stamps the three sheets
# with the new generation timestamps. The Status column on each sheet also
# got noticeably narrower in the regenerated report, so the two "status
# timestamp" columns are resized to match.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Refreshed generation / handoff timestamps ---
$overview.Range("G2").Value = "2016-08-20 01:02:29"
$zhcn.Range("H2").Value     = "2016-08-20 01:02:25"
$dede.Range("H2").Value     = "2016-08-20 01:02:29"

# --- Narrower "status timestamp" columns in the regenerated report ---
# (ColumnWidth is quantized by the host to 1/6-character steps, so use the
# value whose rounded result lands nearest the target OOXML column width of
# 17.2159881591797 characters.)
$newStatusColWidth = 16.333333333333332

$overview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$overview.Columns.Item(6).ColumnWidth = $newStatusColWidth
$zhcn.Columns.Item(3).ColumnWidth     = $newStatusColWidth
$dede.Columns.Item(3).ColumnWidth     = $newStatusColWidth
